$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.812.25'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '2.088.60'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.54'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.48%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.394'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0785'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("E11").Value = '  +3.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.10'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.44%  '
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.35'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.781'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").Value = '2.080.04'
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").Value = '37.757.98'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("E19").Value = '  -0.76%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.55%  '
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  -0.95%  '
$ws.Range("E31").Value = '  +1.02%  '
$ws.Range("E32").Value = '  +0.79%  '
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.67'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("E35").Value = '  -0.59%  '
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0239'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0978'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.29%  '
$ws.Range("E43").Value = '  -1.17%  '
$ws.Range("E44").Value = '  +5.42%  '
$ws.Range("D45").Value = '1.455.90'
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("E46").Value = '  -1.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.81%  '
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("E50").Value = '  -1.54%  '
$ws.Range("D51").Value = '2.280.16'
$ws.Range("E51").Value = '  +0.05%  '
